$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts existing rows 3..50 down to 4..51
$ws.Rows(3).Insert()

# Fill the new row 3 with the new data record
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 45237
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 100112013
$ws.Range("G3").Value = "Alcachofa"
$ws.Range("H3").Value = "Española"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 180
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("N3").Value = "`$/caja 30 unidades"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 433
$ws.Range("Q3").Value = 30
$ws.Range("R3").Value = "Hortaliza"
